$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values
$ws.Range("B2").Value = 0.45355212704528869
$ws.Range("C2").Value = 1.9639324952430477
$ws.Range("D2").Value = 0.41317376586636378
$ws.Range("E2").Value = 0.90333765667704258

# Row 3 data values
$ws.Range("B3").Value = 0.95181260226671338
$ws.Range("C3").Value = 0.75663539154710158
$ws.Range("D3").Value = 0.58953414315802211
$ws.Range("E3").Value = 0.74633447383251594

# Update the selection to match the new highlighted range
$ws.Range("B1:E3").Select()
